# Update the "Förändrad" (changed) date column (C) for rows 2-13
# from 2023-10-08 (serial 45207) to 2023-10-09 (serial 45208).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
